# Updated symbol list on Wed Jan 25 17:59:19 UTC 2023 with GitHub Actions
#
# The sheet stores every data cell as plain text (inline strings in the
# original file: prices like "301.71" and deltas like "-2.60%" are text,
# not numbers). Excel's COM layer auto-detects numeric/percent-looking
# input and would otherwise coerce it to a real number/percentage, so each
# such write is done with a leading apostrophe (forces text entry) and the
# cell style is reset back to "Normal" afterwards so no stray
# quote-prefix/number-format style sticks around.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    $r = $ws.Range($range)
    $r.Value = "'" + $text
    $r.Style = "Normal"
}

# Row 2 - BNB
Set-TextValue "D2" "301.70"
Set-TextValue "E2" "-2.66%"

# Row 3 - OKB
Set-TextValue "D3" "35.35"
Set-TextValue "E3" "-0.63%"

# Row 4 - HuobiToken
Set-TextValue "D4" "5.068"
Set-TextValue "E4" "-0.71%"

# Row 5 - Cronos
Set-TextValue "D5" "0.07933"
Set-TextValue "E5" "-2.78%"

# Row 6
Set-TextValue "D6" "1.882"
Set-TextValue "E6" "-8.25%"

# Row 7
Set-TextValue "D7" "7.785"
Set-TextValue "E7" "-2.00%"

# Row 8
Set-TextValue "D8" "0.9286"
Set-TextValue "E8" "0.37%"

# Row 9
Set-TextValue "D9" "0.1408"
Set-TextValue "E9" "28.50%"

# Row 10
Set-TextValue "D10" "0.1902"
Set-TextValue "E10" "-0.41%"

# Row 11
Set-TextValue "D11" "0.09043"
Set-TextValue "E11" "-1.40%"

# Row 12
Set-TextValue "D12" "0.03437"
Set-TextValue "E12" "-5.85%"

# Row 13
Set-TextValue "D13" "0.09842"
Set-TextValue "E13" "-0.66%"

# Row 14
Set-TextValue "D14" "0.001399"
Set-TextValue "E14" "-2.48%"

# Row 15
Set-TextValue "D15" "0.005887"
Set-TextValue "E15" "2.57%"

# Row 16
Set-TextValue "D16" "3.533"
Set-TextValue "E16" "1.78%"

# Row 17
Set-TextValue "D17" "4.047"
Set-TextValue "E17" "-1.94%"

# Row 18 (only E changes)
Set-TextValue "E18" "4.56%"

# Row 20 (only D, E change; row 19 untouched)
Set-TextValue "D20" "0.1303"
Set-TextValue "E20" "0.20%"

# Row 21
Set-TextValue "D21" "5.049"
Set-TextValue "E21" "-0.88%"

# Row 22
Set-TextValue "D22" "0.2401"
Set-TextValue "E22" "8.47%"

# Row 23
Set-TextValue "D23" "0.04507"
Set-TextValue "E23" "-0.80%"

# Row 24
Set-TextValue "D24" "0.001215"
Set-TextValue "E24" "-0.91%"

# Row 25
Set-TextValue "D25" "0.004759"
Set-TextValue "E25" "-0.50%"

# Row 27 (only E changes)
Set-TextValue "E27" "-32.53%"

# Row 39 (only E changes)
Set-TextValue "E39" "-5.66%"

# Row 40
Set-TextValue "D40" "0.04740"
Set-TextValue "E40" "-2.98%"

# Row 41: was KickToken -> now Dexo
$ws.Range("B41").Value = "Dexo"
$ws.Range("C41").Value = "https://coinranking.com/coin/QkL_pl546+dexo-dexo"
Set-TextValue "D41" "0.01046"
Set-TextValue "E41" "8.01%"

# Row 42: was Dexo -> now KickToken
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D42" "0.007342"
Set-TextValue "E42" "-3.06%"

# Row 43
Set-TextValue "D43" "0.1324"
Set-TextValue "E43" "-4.18%"

# Row 44
Set-TextValue "D44" "0.002112"
Set-TextValue "E44" "-3.93%"

# Row 45
Set-TextValue "D45" "0.01098"
Set-TextValue "E45" "-5.34%"

# Row 46
Set-TextValue "D46" "0.00006226"
Set-TextValue "E46" "-5.53%"

# Row 47 (only E changes)
Set-TextValue "E47" "0.01%"

# Row 48
Set-TextValue "D48" "64.67"
Set-TextValue "E48" "7.33%"

# Row 49 (only E changes)
Set-TextValue "E49" "10.54%"

# Row 50 (only E changes)
Set-TextValue "E50" "0.01%"

# Row 51 (only E changes)
Set-TextValue "E51" "0.01%"
